$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 1559.4762
$ws.Range("I48").Value = 985
$ws.Range("K48").Value = 2955
$ws.Range("M48").Value = -2663

$ws.Range("H56").Value = 1559.4762
$ws.Range("I56").Value = 985
$ws.Range("K56").Value = 2955
$ws.Range("M56").Value = -2421

$ws.Range("H86").Value = 3686
$ws.Range("I86").Value = 3447.75
$ws.Range("K86").Value = 3447.75
$ws.Range("M86").Value = -2324.75

$ws.Range("H89").Value = 3686
$ws.Range("I89").Value = 3447.75
$ws.Range("K89").Value = 17238.75
$ws.Range("M89").Value = -11622.75

$ws.Range("H100").Value = 3054.25
$ws.Range("I100").Value = 2776.4285
$ws.Range("K100").Value = 2776.4285
$ws.Range("M100").Value = -2235.4285

$ws.Range("H106").Value = 10849.167
$ws.Range("I106").Value = 5686.3335
$ws.Range("K106").Value = 5686.3335
$ws.Range("M106").Value = -5055.3335

$ws.Range("H132").Value = 7415.8823
$ws.Range("I132").Value = 8175.3447
$ws.Range("J132").Value = 3011
$ws.Range("K132").Value = 24526.0341
$ws.Range("L132").Value = 9033
$ws.Range("M132").Value = -21996.0341
$ws.Range("N132").Value = -14093

$ws.Range("H135").Value = 400.8889
$ws.Range("I135").Value = 338.625
$ws.Range("J135").Value = 899
$ws.Range("K135").Value = 3047.625
$ws.Range("L135").Value = 8091
$ws.Range("M135").Value = -512.625
$ws.Range("N135").Value = -13161

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19812.717
$ws.Range("I32").Value = 20203.182
$ws.Range("K32").Value = 20203.182
$ws.Range("M32").Value = -19916.182

$ws.Range("H44").Value = 34274.5
$ws.Range("I44").Value = 19000
$ws.Range("J44").Value = 49549
$ws.Range("K44").Value = 19000
$ws.Range("L44").Value = 49549
$ws.Range("M44").Value = -18512
$ws.Range("N44").Value = -50525

$ws.Range("H55").Value = 19499.5
$ws.Range("I55").Value = 19499.5
$ws.Range("K55").Value = 19499.5
$ws.Range("M55").Value = -19184.5

$ws.Range("H61").Value = 6291.5386
$ws.Range("I61").Value = 3572.5454
$ws.Range("K61").Value = 3572.5454
$ws.Range("M61").Value = -3360.5454

$ws.Range("H74").Value = 1041.1666
$ws.Range("I74").Value = 1041.1666
$ws.Range("K74").Value = 1041.1666
$ws.Range("M74").Value = -167.1666

$ws.Range("H77").Value = 1041.1666
$ws.Range("I77").Value = 1041.1666
$ws.Range("K77").Value = 5205.833000000001
$ws.Range("M77").Value = -837.8330000000005

$ws.Range("H80").Value = 86264.336
$ws.Range("J80").Value = 89397
$ws.Range("L80").Value = 89397
$ws.Range("N80").Value = -91393

$ws.Range("H83").Value = 86264.336
$ws.Range("J83").Value = 89397
$ws.Range("L83").Value = 268191
$ws.Range("N83").Value = -278175

$ws.Range("H122").Value = 2903.6667
$ws.Range("I122").Value = 2360.4443
$ws.Range("K122").Value = 7081.3329
$ws.Range("M122").Value = -4631.3329

$ws.Range("H136").Value = 6291.5386
$ws.Range("I136").Value = 3572.5454
$ws.Range("K136").Value = 10717.6362
$ws.Range("M136").Value = -8167.636200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 91977.63
$ws.Range("J22").Value = 1194
$ws.Range("L22").Value = 1194
$ws.Range("N22").Value = -1540

$ws.Range("H80").Value = 494.84616
$ws.Range("J80").Value = 396.3
$ws.Range("L80").Value = 396.3
$ws.Range("N80").Value = -2392.3

$ws.Range("H83").Value = 494.84616
$ws.Range("J83").Value = 396.3
$ws.Range("L83").Value = 1981.5
$ws.Range("N83").Value = -11965.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 12343.333
$ws.Range("J21").Value = 12343.333
$ws.Range("L21").Value = 12343.333
$ws.Range("N21").Value = -12813.333

$ws.Range("H58").Value = 68621.734
$ws.Range("I58").Value = 92535.45
$ws.Range("J58").Value = 2859
$ws.Range("K58").Value = 92535.45
$ws.Range("L58").Value = 2859
$ws.Range("M58").Value = -92332.45
$ws.Range("N58").Value = -3265

$ws.Range("H60").Value = 29343.223
$ws.Range("I60").Value = 4696.6665
$ws.Range("J60").Value = 41666.5
$ws.Range("K60").Value = 4696.6665
$ws.Range("L60").Value = 41666.5
$ws.Range("M60").Value = -4185.6665
$ws.Range("N60").Value = -42688.5

$ws.Range("H107").Value = 1974.1936
$ws.Range("I107").Value = 518.7059
$ws.Range("K107").Value = 518.7059
$ws.Range("M107").Value = 1401.2941

$ws.Range("H112").Value = 79999.5
$ws.Range("J112").Value = 79999.5
$ws.Range("L112").Value = 79999.5
$ws.Range("N112").Value = -82953.5

$ws.Range("H136").Value = 68621.734
$ws.Range("I136").Value = 92535.45
$ws.Range("J136").Value = 2859
$ws.Range("K136").Value = 277606.35
$ws.Range("L136").Value = 8577
$ws.Range("M136").Value = -275056.35
$ws.Range("N136").Value = -13677

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 567.9048
$ws.Range("I113").Value = 719.8333
$ws.Range("J113").Value = 507.13333
$ws.Range("K113").Value = 2159.4999
$ws.Range("L113").Value = 1521.39999
$ws.Range("M113").Value = 10.5001000000002
$ws.Range("N113").Value = -5861.39999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 10833.333
$ws.Range("J10").Value = 16000
$ws.Range("L10").Value = 16000
$ws.Range("N10").Value = -16338

$ws.Range("H22").Value = 8377
$ws.Range("I22").Value = 3008
$ws.Range("J22").Value = 10166.667
$ws.Range("K22").Value = 3008
$ws.Range("L22").Value = 10166.667
$ws.Range("M22").Value = -2479
$ws.Range("N22").Value = -11224.667

$ws.Range("H25").Value = 9833.333000000001
$ws.Range("J25").Value = 9833.333000000001
$ws.Range("L25").Value = 9833.333000000001
$ws.Range("N25").Value = -10891.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2042.4642
$ws.Range("I16").Value = 2144.182
$ws.Range("K16").Value = 2144.182
$ws.Range("M16").Value = -1974.182

$ws.Range("H46").Value = 4867.1455
$ws.Range("I46").Value = 23900.2
$ws.Range("J46").Value = 2963.84
$ws.Range("K46").Value = 23900.2
$ws.Range("L46").Value = 2963.84
$ws.Range("M46").Value = -23712.2
$ws.Range("N46").Value = -3339.84

$ws.Range("H55").Value = 410.1
$ws.Range("I55").Value = 242
$ws.Range("K55").Value = 242
$ws.Range("M55").Value = -69

$ws.Range("H68").Value = 4863.3687
$ws.Range("I68").Value = 3192.6365
$ws.Range("J68").Value = 7160.625
$ws.Range("K68").Value = 3192.6365
$ws.Range("L68").Value = 7160.625
$ws.Range("M68").Value = -2443.6365
$ws.Range("N68").Value = -8658.625

$ws.Range("H71").Value = 4863.3687
$ws.Range("I71").Value = 3192.6365
$ws.Range("J71").Value = 7160.625
$ws.Range("K71").Value = 15963.1825
$ws.Range("L71").Value = 35803.125
$ws.Range("M71").Value = -12219.1825
$ws.Range("N71").Value = -43291.125

$ws.Range("H136").Value = 2789.476
$ws.Range("I136").Value = 1948.9445
$ws.Range("J136").Value = 7832.6665
$ws.Range("K136").Value = 5846.833500000001
$ws.Range("L136").Value = 23497.9995
$ws.Range("M136").Value = -3296.833500000001
$ws.Range("N136").Value = -28597.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 27000.2
$ws.Range("J104").Value = 27000.2
$ws.Range("L104").Value = 27000.2
$ws.Range("N104").Value = -33988.2

$ws.Range("H135").Value = 81607.5
$ws.Range("J135").Value = 81607.5
$ws.Range("L135").Value = 81607.5
$ws.Range("N135").Value = -91747.5
